$wb = $excel.ActiveWorkbook

# Rename sheets
$wb.Worksheets.Item("GNG_TO-16504778668283727").Name = "GNG_TO-16509961603928454"
$wb.Worksheets.Item("NB_TO-16504778686756663").Name = "NB_TO-16509961629381833"
$wb.Worksheets.Item("RS_TO-16504778686816652").Name = "RS_TO-16509961629381833"
$wb.Worksheets.Item("TOL_TO-16504778687406635").Name = "TOL_TO-16509961630022187"
$wb.Worksheets.Item("vSAT_TO-16504778688016996").Name = "vSAT_TO-165099616306623"

# Sheet 1: GNG
$ws1 = $wb.Worksheets.Item("GNG_TO-16509961603928454")
$ws1.Range("B2").Value = "go_stims-16509961603608534.csv"
$ws1.Range("B3").Value = "GNG_stims-1650996160376845.csv"
$ws1.Range("B4").Value = "go_stims-1650996160376845.csv"
$ws1.Range("B5").Value = "GNG_stims-16509961603928454.csv"

# Sheet 2: NB
$ws2 = $wb.Worksheets.Item("NB_TO-16509961629381833")
$ws2.Range("B2").Value = "ZB-match_2-16509961607541826.csv"
$ws2.Range("B3").Value = "ZB-match_9-16509961606341898.csv"
$ws2.Range("B4").Value = "TB-1650996162114223.csv"
$ws2.Range("B5").Value = "OB-16509961615942233.csv"
$ws2.Range("B6").Value = "TB-16509961629141808.csv"
$ws2.Range("B7").Value = "OB-1650996160866185.csv"
$ws2.Range("B8").Value = "OB-16509961616182225.csv"
$ws2.Range("B9").Value = "ZB-match_5-16509961604488792.csv"
$ws2.Range("B10").Value = "TB-16509961620742216.csv"

# Sheet 3: RS - only name changed, no data changes

# Sheet 4: TOL
$ws4 = $wb.Worksheets.Item("TOL_TO-16509961630022187")
$ws4.Range("B2").Value = "MM_stims-1650996162970221.csv"
$ws4.Range("B3").Value = "ZM_stims-16509961629461992.csv"
$ws4.Range("B4").Value = "MM_stims-1650996162986222.csv"
$ws4.Range("B5").Value = "ZM_stims-1650996162970221.csv"
$ws4.Range("B6").Value = "MM_stims-16509961630022187.csv"
$ws4.Range("B7").Value = "ZM_stims-1650996162986222.csv"

# Sheet 5: vSAT
$ws5 = $wb.Worksheets.Item("vSAT_TO-165099616306623")
$ws5.Range("B2").Value = "vSAT_stims-1650996163034219.csv"
$ws5.Range("B3").Value = "SAT_stims-16509961630181842.csv"
$ws5.Range("B4").Value = "vSAT_stims-16509961630502174.csv"
$ws5.Range("B5").Value = "SAT_stims-16509961630022187.csv"
